$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashflow")

# Flip the sign of the literal present-value numbers in columns F:BN
# (col 6 .. col 66) across the four "series" rows (2,3,4,5). Downstream
# formula rows (10, 12, etc.) reference these cells and will recompute
# automatically.
for ($col = 6; $col -le 66; $col++) {
    for ($row = 2; $row -le 5; $row++) {
        $cell = $ws.Cells.Item($row, $col)
        $v = $cell.Value2
        $cell.Value = (0 - $v)
    }
}

# Remove the stray formatted-but-empty row 18 that was left over below the
# table.
$ws.Rows("18:18").Delete() | Out-Null

# Leave the selection where the author last left it.
$ws.Range("J23").Select() | Out-Null
